# Update the header label for the Solar Photovoltaics Generation column:
# it was reported in GWh, correct it to read MWh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Solar Photovoltaics Generation (MWh)"

# Move/save the active cell selection to E5 (was G6).
$ws.Range("E5").Select() | Out-Null
